$d = $word.ActiveDocument

# Locate the anchor paragraphs by their (distinctive) text so the script is
# resilient to any paragraph-index differences. NOTE: the replacement range
# below deliberately starts at the *true* (non-empty-paragraph) start of the
# "Bachelors of Science in Engineering" paragraph and re-emits that
# paragraph's text unchanged -- starting a multi-paragraph InsertXML range
# exactly on the boundary of the (empty) bookmark paragraph that follows it
# causes that empty paragraph to be skipped over rather than replaced.
$anchorPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Bachelors of Science in Engineering") {
        $anchorPara = $d.Paragraphs.Item($i)
    }
    if ($t -eq "Our ability to introduce a new standard in residential accommodation to the market is demonstrated in the experience and expertise of our highly skilled and competent executive management team and technical project team leaders.") {
        $endPara = $d.Paragraphs.Item($i)
    }
}

$r = $d.Range($anchorPara.Range.Start, $endPara.Range.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = @"
<w:p $ns><w:r><w:t>Bachelors of Science in Engineering</w:t></w:r></w:p>
<w:p $ns/>
<w:p $ns>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Sheleisha</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> Kassie </w:t></w:r>
</w:p>
<w:p $ns>
  <w:r>
    <w:br/>
    <w:t xml:space="preserve">Consultant in Marketing, </w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:r><w:t>Interior Design</w:t></w:r>
</w:p>
<w:p $ns>
  <w:r><w:t xml:space="preserve">Space Coordination </w:t></w:r>
</w:p>
<w:p $ns/>
<w:p $ns/>
<w:p $ns/>
<w:p $ns/>
<w:p $ns/>
<w:p $ns/>
<w:p $ns/>
<w:p $ns>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Company Profile</w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:r><w:t xml:space="preserve">WISE emerged as the brainchild of Douglas </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Gourzong</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> who conceptualized and spearheaded the vision</w:t></w:r>
  <w:r><w:t xml:space="preserve"> o</w:t></w:r>
  <w:r><w:t>f a highly technical and comprehensively resourced enterprise; to bring to the region its innovative, state-of-the-art signature residential suite of apartments, condominiums, town houses as well as cutting-edge commercial structures and facilities. The Wise Team boasts qualifications and expertise in all specialized requisite branches if engineering covering Electrical, Mechanical, Civil, Smart Systems, Instrumentation, Aquatics, Architecture and Construction/Project management.</w:t></w:r>
</w:p>
<w:p $ns/>
<w:p $ns>
  <w:r><w:t>Our ability to introduce a new standard in residential accommodation to the market is demonstrated in the experience and expertise of our highly skilled and competent executive management team and technical project team leaders.</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@

$r.InsertXML($xml)
